# Update "想去人数" (column F) values on the 展览, 演出, and 全部类型 sheets
# to reflect the latest scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 2961
$ws1.Range("F9").Value  = 436
$ws1.Range("F11").Value = 56
$ws1.Range("F13").Value = 226
$ws1.Range("F14").Value = 4289
$ws1.Range("F15").Value = 4289
$ws1.Range("F16").Value = 96
$ws1.Range("F17").Value = 86
$ws1.Range("F18").Value = 108
$ws1.Range("F20").Value = 193
$ws1.Range("F22").Value = 6497
$ws1.Range("F24").Value = 97
$ws1.Range("F27").Value = 1229
$ws1.Range("F28").Value = 6246
$ws1.Range("F29").Value = 1629
$ws1.Range("F31").Value = 1859
$ws1.Range("F32").Value = 5960
$ws1.Range("F33").Value = 107
$ws1.Range("F36").Value = 83
$ws1.Range("F37").Value = 404
$ws1.Range("F38").Value = 4081
$ws1.Range("F42").Value = 17
$ws1.Range("F50").Value = 21

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 18

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 2961
$ws4.Range("F9").Value  = 436
$ws4.Range("F11").Value = 56
$ws4.Range("F14").Value = 226
$ws4.Range("F15").Value = 4289
$ws4.Range("F16").Value = 4289
$ws4.Range("F17").Value = 96
$ws4.Range("F18").Value = 86
$ws4.Range("F19").Value = 108
$ws4.Range("F21").Value = 193
$ws4.Range("F22").Value = 6497
$ws4.Range("F24").Value = 97
$ws4.Range("F26").Value = 1229
$ws4.Range("F28").Value = 6246
$ws4.Range("F29").Value = 1629
$ws4.Range("F32").Value = 1859
$ws4.Range("F33").Value = 5960
$ws4.Range("F34").Value = 107
$ws4.Range("F37").Value = 83
$ws4.Range("F38").Value = 404
$ws4.Range("F39").Value = 4081
$ws4.Range("F42").Value = 17
$ws4.Range("F50").Value = 18
